$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each row with refreshed crypto data
$ws.Range("D2").Value = "27.898.77"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "1.783.38"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'310.67"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.5117"
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("D8").Value = "'0.3765"
$ws.Range("E8").Value = "  -2.23%  "
$ws.Range("D9").Value = "'0.07774"
$ws.Range("E9").Value = "  -8.49%  "
$ws.Range("D10").Value = "'41.34"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("D11").Value = "'1.084"
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'6.184"
$ws.Range("E13").Value = "  -3.80%  "
$ws.Range("D14").Value = "'20.17"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").Value = "1.779.97"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "'7.182"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").Value = "'91.94"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("E18").Value = "  -6.17%  "
$ws.Range("D19").Value = "'0.06530"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "'17.01"
$ws.Range("E21").Value = "  -4.01%  "
$ws.Range("D22").Value = "'5.892"
$ws.Range("E22").Value = "  -3.07%  "
$ws.Range("D23").Value = "27.945.48"
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").Value = "'10.94"
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").Value = "'2.245"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").Value = "'158.21"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'20.23"
$ws.Range("E27").Value = "  -4.73%  "
$ws.Range("D28").Value = "1.980.10"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").Value = "'2.349"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "'122.19"
$ws.Range("D31").Value = "'0.1071"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").Value = "'1.039"
$ws.Range("E32").Value = "  -5.09%  "
$ws.Range("D33").Value = "'3.631"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "'5.476"
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("D35").Value = "'0.07087"
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("D36").Value = "'0.02304"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "'0.2122"
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("D38").Value = "'8.575"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D41").Value = "'0.6096"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("D42").Value = "'1.153"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("D44").Value = "'0.5962"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("D45").Value = "'13.07"
$ws.Range("E45").Value = "  -3.32%  "
$ws.Range("D46").Value = "'3.727"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").Value = "'126.37"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").Value = "'1.213"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "'1.893"
$ws.Range("E49").Value = "  -4.75%  "
$ws.Range("D50").Value = "'0.06724"
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("D51").Value = "'1.053"
$ws.Range("E51").Value = "  -1.74%  "

# Rows 39/40: source data coin order changed (Aptos <-> InternetComputer(DFINITY))
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'11.49"
$ws.Range("E39").Value = "  +2.17%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'5.006"
$ws.Range("E40").Value = "  -4.00%  "
